# Added comments, improved interaction: append four new trade rows
# (three AAPL entries and one A entry) below the existing trade log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: TICKER, DATE, BUY/SELL, PRICE, VOLUME, NET_EFFECT_TO_CASH,
# TOTAL_SHARES_HOLDING, TICKER_TOTAL_VALUE, AVERAGE_PRICE
$newTrades = @(
    @{ Ticker = "AAPL"; Date = "08/15/20"; Side = "SELL"; Price = 200;  Volume = 20;  NetCash = 4000;   Shares = 20;  TotalValue = 4000;   AvgPrice = 200 },
    @{ Ticker = "AAPL"; Date = "12/12/12"; Side = "BUY";  Price = 200;  Volume = 200; NetCash = -40000; Shares = 220; TotalValue = 44000;  AvgPrice = 200 },
    @{ Ticker = "AAPL"; Date = "08/15/20"; Side = "BUY";  Price = 12;   Volume = 13;  NetCash = -156;   Shares = 233; TotalValue = 44156;  AvgPrice = 189.51 },
    @{ Ticker = "A";    Date = "08/15/20"; Side = "BUY";  Price = 12;   Volume = 12;  NetCash = -144;   Shares = 12;  TotalValue = 144;    AvgPrice = 12 }
)

$startRow = 8
$endRow = $startRow + $newTrades.Count - 1

# The DATE column stores plain "MM/DD/YY" text (like the rest of the
# sheet), not real date serials. Briefly force the new date cells to
# text so Excel's autodetection doesn't reinterpret them as dates, then
# clear the format again so no stray formatting sticks around.
$dateRange = $ws.Range("B$startRow`:B$endRow")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $newTrades.Count; $i++) {
    $row = $startRow + $i
    $trade = $newTrades[$i]

    $ws.Cells.Item($row, 1).Value = $trade.Ticker
    $ws.Cells.Item($row, 2).Value = $trade.Date
    $ws.Cells.Item($row, 3).Value = $trade.Side
    $ws.Cells.Item($row, 4).Value = $trade.Price
    $ws.Cells.Item($row, 5).Value = $trade.Volume
    $ws.Cells.Item($row, 6).Value = $trade.NetCash
    $ws.Cells.Item($row, 7).Value = $trade.Shares
    $ws.Cells.Item($row, 8).Value = $trade.TotalValue
    $ws.Cells.Item($row, 9).Value = $trade.AvgPrice
    # Column J (REALIZED_PROFIT) is left blank for these rows, same as
    # the rest of the trade log.
}

$dateRange.ClearFormats()
